# feat: add 2022-Q1 data
#
# The workbook currently has 3 sheets: 2021-Q1, 2021-Q2, 总计 (Total).
# We add a new "2022-Q1" sheet (per-fund holding detail, same shape as
# the 2021-Q1 / 2021-Q2 sheets) right before the "总计" summary sheet,
# and refresh "总计" so it aggregates 2022-Q1 + 2021-Q2 + 2021-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet shuffle: rename the existing "总计" sheet to "2022-Q1" (so it
#    keeps its original sheetId / file identity), then clone it right
#    after itself (Copy carries over sheetPr/pageSetUpPr/etc. exactly)
#    and rename the clone back to "总计" to hold the refreshed summary.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Name = "2022-Q1"
$totalSheet.Copy($null, $totalSheet)

$wsQ1_2022 = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("2022-Q1 (2)")
$wsTotal.Name = "总计"

$wsSrc = $wb.Worksheets.Item("2021-Q2")

# ---------------------------------------------------------------------
# 2) "2022-Q1" sheet: wipe the old aggregate-style content it inherited
#    and rebuild it as a per-fund holdings sheet (same column layout as
#    2021-Q1 / 2021-Q2: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#    持有市值(亿元)/仓位排名).
# ---------------------------------------------------------------------
$wsQ1_2022.Cells.Clear()

# Pull header (row1 B:H) and column-A formatting from the 2021-Q2 sheet,
# which already carries the bold/centered/bordered "s=2" style we need.
$wsSrc.Range("B1:H1").Copy()
$wsQ1_2022.Range("B1:H1").PasteSpecial(-4122)
$wsSrc.Range("A2").Copy()
$wsQ1_2022.Range("A2").PasteSpecial(-4122)

$wsQ1_2022.Range("B1").Value = "基金代码"
$wsQ1_2022.Range("C1").Value = "基金名称"
$wsQ1_2022.Range("D1").Value = "基金规模"
$wsQ1_2022.Range("E1").Value = "股票总仓位"
$wsQ1_2022.Range("F1").Value = "仓位占比"
$wsQ1_2022.Range("G1").Value = "持有市值(亿元)"
$wsQ1_2022.Range("H1").Value = "仓位排名"

$wsQ1_2022.Range("A2").Value = 0

# B2/D2/E2/F2/G2 look like numbers (leading-zero fund code, plain
# decimals) -- force them to be stored as text (with no cell style left
# behind) so they round-trip exactly like the source data, matching the
# inline-string cells in the target sheet.
$textRange = $wsQ1_2022.Range("B2:G2")
$textRange.NumberFormat = "@"
$wsQ1_2022.Range("B2").Value = "006157"
$wsQ1_2022.Range("C2").Value = "财通量化核心优选混合"
$wsQ1_2022.Range("D2").Value = "0.09"
$wsQ1_2022.Range("E2").Value = "92.85"
$wsQ1_2022.Range("F2").Value = "1.43"
$wsQ1_2022.Range("G2").Value = "0.0013"
$textRange.ClearFormats()

$wsQ1_2022.Range("H2").Value = 7

# ---------------------------------------------------------------------
# 3) "总计" sheet: rebuild the per-quarter roll-up, now with 2022-Q1 on
#    top, followed by 2021-Q2 and 2021-Q1.
# ---------------------------------------------------------------------
$wsTotal.Cells.Clear()

$wsSrc.Range("B1:D1").Copy()
$wsTotal.Range("B1:D1").PasteSpecial(-4122)
$wsSrc.Range("A2").Copy()
$wsTotal.Range("A2:A4").PasteSpecial(-4122)

$wsTotal.Range("B1").Value = "日期"
$wsTotal.Range("C1").Value = "持有数量(只)"
$wsTotal.Range("D1").Value = "持有市值(亿元)"

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 1
$wsTotal.Range("D2").Value = 0

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2021-Q2"
$wsTotal.Range("C3").Value = 7
$wsTotal.Range("D3").Value = 0.23

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2021-Q1"
$wsTotal.Range("C4").Value = 4
$wsTotal.Range("D4").Value = 0.01

# ---------------------------------------------------------------------
# 4) Restore the originally-active sheet/tab (2021-Q1) as selected, so
#    we don't leave the newly-touched sheets marked as the active tab.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q1").Activate()
